$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 22: add "Y" to copy (C) and delete (D) columns ---
$ws.Range("C22").Value = "Y"
$ws.Range("D22").Value = "Y"

# --- Row 25: add "Y" to copy (C) and delete (D) columns ---
$ws.Range("C25").Value = "Y"
$ws.Range("D25").Value = "Y"

# --- New row 33: a "blank-ish" delete marker (two spaces) ---
# Create this brand-new shared string BEFORE "No Data" so it lands at the
# same shared-string index (175) that the target workbook expects.
$ws.Range("D33").Value = "  "

# --- Rows 36-42: add "Y" job-submission marker ---
$ws.Range("B36").Value = "Y"
$ws.Range("B37").Value = "Y"
$ws.Range("B38").Value = "Y"
$ws.Range("B39").Value = "Y"
$ws.Range("B40").Value = "Y"
$ws.Range("B41").Value = "Y"
$ws.Range("B42").Value = "Y"

# --- Rows 43-44: "No Data" job-submission marker (new shared string, index 176) ---
$ws.Range("B43").Value = "No Data"
$ws.Range("B44").Value = "No Data"

# --- Rows 45-48: add "Y" job-submission marker ---
$ws.Range("B45").Value = "Y"
$ws.Range("B46").Value = "Y"
$ws.Range("B47").Value = "Y"
$ws.Range("B48").Value = "Y"

# --- Update the view: scroll down a bit and move the selection ---
$ws.Range("I29").Select()
